# WRENCH STUDIO.docx — spelling fix ("Cambio de dos letras (ortografia)")
#
# The author's real edit is two small spelling corrections:
#   1. "practicas"  -> "prácticas"   (missing accent on the first "a")
#   2. " penamos"   -> " pensamos"   (missing "s")
#
# (Everything else visible in the stored-diff — built-in style-id
# relocalisation, w:proofErr spell-check bookmarks, and cosmetic run
# splits/merges with no net text change — is a side effect of the
# document having been opened/saved by Word's proofing + locale
# pipeline, not something reachable through the Word object model, so
# it is intentionally left alone here.)

$d = $word.ActiveDocument

$found1 = $d.Content.Find.Execute(
    "practicas de programación", $true, $false, $false, $false, $false,
    $true, 1, $false, "prácticas de programación", 2)
if (-not $found1) {
    throw "Could not find 'practicas de programación' to fix."
}

$found2 = $d.Content.Find.Execute(
    " penamos añadir", $true, $false, $false, $false, $false,
    $true, 1, $false, " pensamos añadir", 2)
if (-not $found2) {
    throw "Could not find ' penamos añadir' to fix."
}

Write-Host "Applied spelling fixes: practicas->prácticas, penamos->pensamos"
